$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Insert a new worksheet "2022-Q1" right after "2021-Q4" (before "总计")
#    and populate it with the fund-holdings detail for that quarter.
# ---------------------------------------------------------------------------
$afterSheet = $wb.Worksheets.Item("2021-Q4")
$newSheet = $wb.Worksheets.Add($null, $afterSheet)
$newSheet.Name = "2022-Q1"

# Match the sheetPr / page setup conventions used by the other quarter sheets.
$newSheet.Outline.SummaryRow = 1
$newSheet.Outline.SummaryColumn = 1
$newSheet.PageSetup.LeftMargin = 54
$newSheet.PageSetup.RightMargin = 54
$newSheet.PageSetup.TopMargin = 72
$newSheet.PageSetup.BottomMargin = 72
$newSheet.PageSetup.HeaderMargin = 36
$newSheet.PageSetup.FooterMargin = 36

# Re-use the header / index-column formatting already present on another
# quarter sheet so the new sheet matches the existing look & feel.
$srcFmt = $wb.Worksheets.Item("2021-Q4")
$srcFmt.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$srcFmt.Range("A2").Copy()
$newSheet.Range("A2:A35").PasteSpecial(-4122)

$headers = @("基金代码","基金名称","基金规模","股票总仓位","仓位占比","持有市值(亿元)","仓位排名")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $newSheet.Cells.Item(1, $c + 2).Value = $headers[$c]
}

$fundRows = @(
    ,@(0,'163407','兴全沪深300指数增强(LOF)A','41.45','95.41','3.72','1.5419',6)
    ,@(1,'159997','天弘中证电子ETF','11.54','99.81','2.50','0.2885',10)
    ,@(2,'515150','富国中证国企一带一路ETF','7.41','99.21','2.59','0.1919',5)
    ,@(3,'515750','富国中证科技50策略ETF','6.60','99.46','2.81','0.1855',8)
    ,@(4,'515860','嘉实中证新兴科技100策略ETF','2.25','98.94','7.67','0.1726',2)
    ,@(5,'515260','华宝中证电子50ETF','4.60','98.50','3.22','0.1481',10)
    ,@(6,'515110','易方达中证国企一带一路ETF','4.83','99.52','2.59','0.1251',5)
    ,@(7,'012124','博道盛彦混合型证券投资基金A','3.40','88.74','3.30','0.1122',9)
    ,@(8,'159786','银华中证虚拟现实主题交易型开放式指数证券投资基金','1.91','97.95','5.81','0.1110',5)
    ,@(9,'010202','天弘中证科技100指数增强A','2.86','94.27','3.76','0.1075',5)
    ,@(10,'010908','大成沪深300指数增强A','3.00','86.52','3.46','0.1038',6)
    ,@(11,'159916','建信深证基本面60ETF','4.09','98.79','2.18','0.0892',10)
    ,@(12,'010203','天弘中证科技100指数增强C','2.17','94.27','3.76','0.0816',5)
    ,@(13,'159910','嘉实深证基本面120ETF','4.67','99.44','1.68','0.0785',10)
    ,@(14,'163116','申万菱信中证申万电子行业投资指数（LOF）A','2.27','93.21','2.75','0.0624',9)
    ,@(15,'012319','西藏东财国证消费电子主题指数增强A','1.12','92.75','5.49','0.0615',5)
    ,@(16,'515320','华安中证电子50ETF','1.88','96.98','3.17','0.0596',10)
    ,@(17,'007230','兴全沪深300指数增强(LOF)C','1.38','95.41','3.72','0.0513',6)
    ,@(18,'159709','工银瑞信深证物联网50交易型开放式指数证券投资基金','1.30','98.81','3.88','0.0504',10)
    ,@(19,'010909','大成沪深300指数增强C','1.41','86.52','3.46','0.0488',6)
    ,@(20,'012320','西藏东财国证消费电子主题指数增强C','0.66','92.75','5.49','0.0362',5)
    ,@(21,'515870','嘉实中证先进制造100策略ETF','0.42','98.79','6.92','0.0291',5)
    ,@(22,'515990','汇添富中证国企一带一路ETF','1.08','99.16','2.61','0.0282',5)
    ,@(23,'000042','中证财通中国可持续发展100 (ECPI ESG) 指数增强A','1.82','94.72','1.32','0.0240',10)
    ,@(24,'517360','华安中证沪港深科技100交易型开放式指数证券投资基金','0.72','93.71','2.84','0.0204',10)
    ,@(25,'003749','创金合信鑫收益灵活配置混合A','0.65','51.22','2.19','0.0142',4)
    ,@(26,'006906','创金合信鑫收益灵活配置混合E','0.65','51.22','2.19','0.0142',4)
    ,@(27,'159913','交银深证300价值ETF','0.52','97.62','2.20','0.0114',10)
    ,@(28,'005035','银华信息科技量化优选股票A','0.26','90.98','2.39','0.0062',4)
    ,@(29,'012125','博道盛彦混合型证券投资基金C','0.15','88.74','3.30','0.0050',9)
    ,@(30,'010531','申万菱信中证申万电子行业投资指数（LOF）C','0.14','93.21','2.75','0.0038',9)
    ,@(31,'005036','银华信息科技量化优选股票C','0.04','90.98','2.39','0.0010',4)
    ,@(32,'003750','创金合信鑫收益灵活配置混合C','0.02','51.22','2.19','0.0004',4)
    ,@(33,'003184','中证财通中国可持续发展100 (ECPI ESG) 指数增强C','0.00','94.72','1.32','0',10)
)

# Columns B:G hold codes/figures that look numeric but are stored as text in
# the source data, so force a text format before assigning the values.
$newSheet.Range("B2:G35").NumberFormat = "@"

for ($i = 0; $i -lt $fundRows.Length; $i++) {
    $row = $fundRows[$i]
    $r = $i + 2
    $newSheet.Cells.Item($r, 1).Value = $row[0]
    $newSheet.Cells.Item($r, 2).Value = $row[1]
    $newSheet.Cells.Item($r, 3).Value = $row[2]
    $newSheet.Cells.Item($r, 4).Value = $row[3]
    $newSheet.Cells.Item($r, 5).Value = $row[4]
    $newSheet.Cells.Item($r, 6).Value = $row[5]
    $newSheet.Cells.Item($r, 7).Value = $row[6]
    $newSheet.Cells.Item($r, 8).Value = $row[7]
}

# Drop the temporary text-number-format style again (keep the text values).
$newSheet.Range("B2:G35").Style = "Normal"

# The very last fund's holding value is exactly 0, stored as a real number.
$newSheet.Range("G35").Value = 0

# ---------------------------------------------------------------------------
# 2) Add the 2022-Q1 summary row to the "总计" sheet, pushing the older
#    quarters down by one row.
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")

# Make sure the new bottom row (row 7) gets the same "index column" style
# that the existing rows already use, before we write any values into it.
$totalSheet.Range("A2").Copy()
$totalSheet.Range("A7").PasteSpecial(-4122)

$summaryRows = @(
    ,@(0,"2022-Q1",34,3.87)
    ,@(1,"2021-Q4",20,3.46)
    ,@(2,"2021-Q3",49,8.27)
    ,@(3,"2021-Q2",70,12.57)
    ,@(4,"2021-Q1",150,60.68)
    ,@(5,"2020-Q4",65,15.64)
)

for ($i = 0; $i -lt $summaryRows.Length; $i++) {
    $row = $summaryRows[$i]
    $r = $i + 2
    $totalSheet.Cells.Item($r, 1).Value = $row[0]
    $totalSheet.Cells.Item($r, 2).Value = $row[1]
    $totalSheet.Cells.Item($r, 3).Value = $row[2]
    $totalSheet.Cells.Item($r, 4).Value = $row[3]
}

# ---------------------------------------------------------------------------
# 3) Restore the originally active sheet/tab.
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("2020-Q4").Activate()

Write-Output "2022-Q1 data added"
